$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I10").Value = 'sd'
$ws.Range("J10").Value = 'Statement-non-opinion'
$ws.Range("I19").Value = 'sv'
$ws.Range("J19").Value = 'Statement-opinion'
$ws.Range("I41").Value = 'sv'
$ws.Range("J41").Value = 'Statement-opinion'
$ws.Range("I45").Value = 'sd'
$ws.Range("J45").Value = 'Statement-non-opinion'
$ws.Range("I46").Value = 'sd'
$ws.Range("J46").Value = 'Statement-non-opinion'
$ws.Range("I59").Value = 'sd'
$ws.Range("J59").Value = 'Statement-non-opinion'
$ws.Range("I66").Value = 'sv'
$ws.Range("J66").Value = 'Statement-opinion'
$ws.Range("I100").Value = 'sd'
$ws.Range("J100").Value = 'Statement-non-opinion'
$ws.Range("I103").Value = 'ba'
$ws.Range("J103").Value = 'Appreciation'
$ws.Range("I105").Value = 'sv'
$ws.Range("J105").Value = 'Statement-opinion'
$ws.Range("I111").Value = 'sv'
$ws.Range("J111").Value = 'Statement-opinion'
$ws.Range("I115").Value = 'sv'
$ws.Range("J115").Value = 'Statement-opinion'
$ws.Range("I122").Value = 'sv'
$ws.Range("J122").Value = 'Statement-opinion'
$ws.Range("I123").Value = 'sv'
$ws.Range("J123").Value = 'Statement-opinion'
$ws.Range("I134").Value = 'sv'
$ws.Range("J134").Value = 'Statement-opinion'
$ws.Range("I165").Value = 'sv'
$ws.Range("J165").Value = 'Statement-opinion'
$ws.Range("I166").Value = 'sv'
$ws.Range("J166").Value = 'Statement-opinion'
$ws.Range("I167").Value = 'sv'
$ws.Range("J167").Value = 'Statement-opinion'
$ws.Range("I171").Value = 'sd'
$ws.Range("J171").Value = 'Statement-non-opinion'
$ws.Range("I173").Value = 'sd'
$ws.Range("J173").Value = 'Statement-non-opinion'
$ws.Range("I180").Value = 'b'
$ws.Range("J180").Value = 'Acknowledge (Backchannel)'
$ws.Range("I196").Value = 'sv'
$ws.Range("J196").Value = 'Statement-opinion'
$ws.Range("I200").Value = 'sv'
$ws.Range("J200").Value = 'Statement-opinion'
$ws.Range("I214").Value = 'qy'
$ws.Range("J214").Value = 'Yes-No-Question'
$ws.Range("I226").Value = 'sv'
$ws.Range("J226").Value = 'Statement-opinion'
$ws.Range("I229").Value = 'sv'
$ws.Range("J229").Value = 'Statement-opinion'
$ws.Range("I232").Value = 'sd'
$ws.Range("J232").Value = 'Statement-non-opinion'
$ws.Range("I237").Value = 'sv'
$ws.Range("J237").Value = 'Statement-opinion'
$ws.Range("I257").Value = 'sv'
$ws.Range("J257").Value = 'Statement-opinion'
$ws.Range("I261").Value = 'sv'
$ws.Range("J261").Value = 'Statement-opinion'
$ws.Range("I265").Value = 'aa'
$ws.Range("J265").Value = 'Agree/Accept'
$ws.Range("I266").Value = 'sv'
$ws.Range("J266").Value = 'Statement-opinion'
$ws.Range("I276").Value = 'sd'
$ws.Range("J276").Value = 'Statement-non-opinion'
$ws.Range("I306").Value = 'sv'
$ws.Range("J306").Value = 'Statement-opinion'
$ws.Range("I310").Value = 'ba'
$ws.Range("J310").Value = 'Appreciation'
$ws.Range("I315").Value = 'sd'
$ws.Range("J315").Value = 'Statement-non-opinion'
$ws.Range("I316").Value = 'ba'
$ws.Range("J316").Value = 'Appreciation'
$ws.Range("I339").Value = 'aa'
$ws.Range("J339").Value = 'Agree/Accept'
$ws.Range("I340").Value = 'aa'
$ws.Range("J340").Value = 'Agree/Accept'
$ws.Range("I343").Value = 'sv'
$ws.Range("J343").Value = 'Statement-opinion'
$ws.Range("I358").Value = '%'
$ws.Range("J358").Value = 'Uninterpretable'
$ws.Range("I385").Value = 'aa'
$ws.Range("J385").Value = 'Agree/Accept'
$ws.Range("I403").Value = 'sv'
$ws.Range("J403").Value = 'Statement-opinion'
$ws.Range("I405").Value = 'sd'
$ws.Range("J405").Value = 'Statement-non-opinion'
$ws.Range("I407").Value = 'sv'
$ws.Range("J407").Value = 'Statement-opinion'
$ws.Range("I408").Value = 'aa'
$ws.Range("J408").Value = 'Agree/Accept'
$ws.Range("I409").Value = 'sd'
$ws.Range("J409").Value = 'Statement-non-opinion'
$ws.Range("I412").Value = 'sv'
$ws.Range("J412").Value = 'Statement-opinion'
$ws.Range("I420").Value = 'qy'
$ws.Range("J420").Value = 'Yes-No-Question'
$ws.Range("I423").Value = 'sv'
$ws.Range("J423").Value = 'Statement-opinion'
$ws.Range("I424").Value = 'ba'
$ws.Range("J424").Value = 'Appreciation'
$ws.Range("I426").Value = 'b'
$ws.Range("J426").Value = 'Acknowledge (Backchannel)'
$ws.Range("I428").Value = 'sd'
$ws.Range("J428").Value = 'Statement-non-opinion'
$ws.Range("I429").Value = 'sd'
$ws.Range("J429").Value = 'Statement-non-opinion'
$ws.Range("I439").Value = 'sv'
$ws.Range("J439").Value = 'Statement-opinion'
$ws.Range("I454").Value = 'sv'
$ws.Range("J454").Value = 'Statement-opinion'
$ws.Range("I457").Value = 'sv'
$ws.Range("J457").Value = 'Statement-opinion'
$ws.Range("I464").Value = 'sd'
$ws.Range("J464").Value = 'Statement-non-opinion'
$ws.Range("I472").Value = 'aa'
$ws.Range("J472").Value = 'Agree/Accept'
$ws.Range("I475").Value = 'ba'
$ws.Range("J475").Value = 'Appreciation'
$ws.Range("I477").Value = 'ba'
$ws.Range("J477").Value = 'Appreciation'
$ws.Range("I479").Value = 'sv'
$ws.Range("J479").Value = 'Statement-opinion'
$ws.Range("I481").Value = 'sv'
$ws.Range("J481").Value = 'Statement-opinion'
$ws.Range("I484").Value = 'sv'
$ws.Range("J484").Value = 'Statement-opinion'
$ws.Range("I504").Value = 'sd'
$ws.Range("J504").Value = 'Statement-non-opinion'
$ws.Range("I506").Value = 'sv'
$ws.Range("J506").Value = 'Statement-opinion'
$ws.Range("I511").Value = 'sd'
$ws.Range("J511").Value = 'Statement-non-opinion'
$ws.Range("I539").Value = 'b'
$ws.Range("J539").Value = 'Acknowledge (Backchannel)'
